# Fix Lines of Code values for human_written rows (rows 23-29, column C)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C23").Value = 890
$ws.Range("C24").Value = 1626
$ws.Range("C25").Value = 163
$ws.Range("C26").Value = 3301
$ws.Range("C27").Value = 3308
$ws.Range("C28").Value = 432
$ws.Range("C29").Value = 745
